$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'5.55%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'32.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'9.75%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.259"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.24%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07461"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'6.88%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.852"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'5.58%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.790"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'6.72%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.534"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'9.09%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'2.00%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.01744"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2,590.46%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1680"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'4.30%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08027"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'6.45%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.07960"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.48%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03025"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'2.99%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09897"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'9.86%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001497"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-6.07%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.04605"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'2.02%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006224"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.06%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'0.07%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'2.231"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.11%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'2.62%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'0.73%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.484"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'11.52%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.1621"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.38%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001219"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.01%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004451"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'4.88%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'19.99%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001748"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'4.96%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04491"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.95%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007166"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.38%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1349"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'8.25%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002191"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'6.10%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01282"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'10.97%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006149"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'5.63%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.7093"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-63.24%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.01299"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.33%"
$ws.Range("E47").Style = "Normal"
